# Weekly fruit/vegetable price update: swap the Fecha/Volumen/Unidad de
# comercializacion/Precio $ por Kg/Kg por unidad values between several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($addr1, $addr2) {
    $r1 = $ws.Range($addr1)
    $r2 = $ws.Range($addr2)
    $tmp = $r1.Value2
    $r1.Value = $r2.Value2
    $r2.Value = $tmp
}

# Rows 2 and 5: swap Fecha (D) and Volumen (M)
Swap-Cell "D2" "D5"
Swap-Cell "M2" "M5"

# Rows 3 and 9: swap Fecha (D), Volumen (M), Unidad de comercializacion (Q),
# Precio $/Kg (S) and Kg / unidad (T)
Swap-Cell "D3" "D9"
Swap-Cell "M3" "M9"
Swap-Cell "Q3" "Q9"
Swap-Cell "S3" "S9"
Swap-Cell "T3" "T9"

# Rows 8 and 10: swap Fecha (D) and Volumen (M)
Swap-Cell "D8" "D10"
Swap-Cell "M8" "M10"
